# Apply the dated-worksheet refresh: update the header date and the
# twenty-four division problems in the table.
#
# NOTE: one pair of replacements has a collision - the text produced by
# one replacement ("33÷3=") is also the *source* text of a later
# replacement elsewhere in the document. To avoid the second find
# matching text we just wrote, that pair is executed first.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Must run before "95÷9=" -> "33÷3=" (see note above).
Replace-Text "33÷3=" "34÷6="

Replace-Text "2025-05-17 Saturday" "2025-05-18 Sunday"

Replace-Text "50÷9=" "22÷7="
Replace-Text "62÷5=" "42÷2="
Replace-Text "16÷6=" "22÷4="
Replace-Text "24÷6=" "86÷5="
Replace-Text "71÷7=" "35÷5="

Replace-Text "24÷8=" "15÷4="
Replace-Text "68÷8=" "35÷6="
Replace-Text "13÷5=" "29÷8="
Replace-Text "95÷9=" "33÷3="
Replace-Text "72÷5=" "83÷4="

Replace-Text "61÷2=" "60÷7="
Replace-Text "73÷6=" "10÷4="
Replace-Text "76÷7=" "99÷8="
Replace-Text "10÷5=" "11÷4="
Replace-Text "42÷3=" "68÷9="

Replace-Text "83÷3=" "23÷8="
Replace-Text "75÷5=" "13÷4="
Replace-Text "82÷9=" "12÷2="
Replace-Text "78÷9=" "28÷5="
Replace-Text "91÷8=" "68÷9="

Replace-Text "38÷8=" "47÷7="
Replace-Text "37÷7=" "43÷9="
Replace-Text "84÷8=" "40÷8="
Replace-Text "58÷3=" "59÷2="
